$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 18: new Single Moth entry
$ws.Range("A18").Value = "2025_01_30"
$ws.Range("B18").Value = "Male"
$ws.Range("C18").Value = 2.12
$ws.Range("D18").Value = 4.44
$ws.Range("E18").Value = 3.838
$ws.Range("H18").Value = 2
$ws.Range("J18").Value = 1
$ws.Range("K18").Value = 1
$ws.Range("L18").Value = 2
$ws.Range("M18").Value = 10000
$ws.Range("N18").Value = 125000
